# "Generate Report for Handback"
#
# Updates the localization-status report after a handback:
#   - Status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it appears
#     (Overview sheet's per-language status column, and each
#     language sheet's "Status" column).
#   - Each language sheet gains two new populated columns,
#     "Latest Target File" (F) and "Latest Handback File" (G),
#     with hyperlinks to the handed-back source + translation
#     files, for both data rows.
#   - "Latest Handback DateTime" (H) is stamped with the handback
#     timestamp for each language sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$newStatus = "Handed back: in sync with en-US"

# ---- Status column: "Ready for handoff" -> "Handed back: in sync with en-US" ----
$ws1.Range("B2").Value = $newStatus
$ws1.Range("C2").Value = $newStatus
$ws1.Range("B3").Value = $newStatus
$ws1.Range("C3").Value = $newStatus

$ws2.Range("C2").Value = $newStatus
$ws2.Range("C3").Value = $newStatus

$ws3.Range("C2").Value = $newStatus
$ws3.Range("C3").Value = $newStatus

# ---- Latest Handback DateTime (column H) ----
$ws2.Range("H2").Value = "2016-03-19 04:29:00"
$ws2.Range("H3").Value = "2016-03-19 04:29:00"

$ws3.Range("H2").Value = "2016-03-19 04:29:06"
$ws3.Range("H3").Value = "2016-03-19 04:29:06"

# ---- New "Latest Target File" (F) / "Latest Handback File" (G) columns ----

# zh-cn
$zhMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/a7bee257839d1b52d65c5101ee59d3b60e9ba40b/e2e/a.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/600b74cf7db7778940aa6c59d572289dd7f20ef7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$ws2.Hyperlinks.Add($ws2.Range("F2"), $zhMdUrl, "", "", "a.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), $zhXlfUrl, "", "", $zhXlfName)
$ws2.Hyperlinks.Add($ws2.Range("F3"), $zhMdUrl, "", "", "a.md")
$ws2.Hyperlinks.Add($ws2.Range("G3"), $zhXlfUrl, "", "", $zhXlfName)

# de-de
$deMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/a7bee257839d1b52d65c5101ee59d3b60e9ba40b/e2e/a.md"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b76ee44d935cae272739f1587cc0df8141d4a33a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$ws3.Hyperlinks.Add($ws3.Range("F2"), $deMdUrl, "", "", "a.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), $deXlfUrl, "", "", $deXlfName)
$ws3.Hyperlinks.Add($ws3.Range("F3"), $deMdUrl, "", "", "a.md")
$ws3.Hyperlinks.Add($ws3.Range("G3"), $deXlfUrl, "", "", $deXlfName)

Write-Output "Handback report generated"
